# Remove the 2007, 2008, 2009 data rows (rows 2-4), shifting the
# remaining 2010-2013 rows up so the table becomes A1:F5 instead of A1:F8.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:A4").EntireRow.Delete()
